$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.260.50"
$ws.Range("E2").Value = "  +3.74%  "

$ws.Range("D3").Value = "1.590.81"
$ws.Range("E3").Value = "  +1.72%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.71"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("E6").Value = "  +1.05%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.17"
$ws.Range("E8").Value = "  +8.84%  "

$ws.Range("E9").Value = "  +0.68%  "

$ws.Range("E10").Value = "  +0.97%  "

$ws.Range("E11").Value = "  +2.00%  "

$ws.Range("E12").Value = "  +1.76%  "

$ws.Range("D13").Value = "1.578.85"
$ws.Range("E13").Value = "  +0.71%  "

$ws.Range("E14").Value = "  +2.40%  "

$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "28.312.96"
$ws.Range("E16").Value = "  +3.89%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.19"
$ws.Range("E17").Value = "  +2.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "227.97"
$ws.Range("E18").Value = "  +4.50%  "

$ws.Range("D19").Value = "0.0₃0710"
$ws.Range("E19").Value = "  +1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.49"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("E21").Value = "  -0.07%  "

$ws.Range("E22").Value = "  -0.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.34"
$ws.Range("E23").Value = "  -0.49%  "

$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.89"
$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.22"
$ws.Range("E26").Value = "  +1.31%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.108"
$ws.Range("E27").Value = "  +0.85%  "

$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("E29").Value = "  -0.05%  "

$ws.Range("E30").Value = "  -0.49%  "

$ws.Range("E31").Value = "  +0.75%  "

$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("D34").Value = "1.403.88"
$ws.Range("E34").Value = "  -3.74%  "

$ws.Range("E35").Value = "  -1.91%  "

$ws.Range("E36").Value = "  -6.73%  "

$ws.Range("E37").Value = "  +0.10%  "

$ws.Range("E38").Value = "  +0.60%  "

$ws.Range("E39").Value = "  +8.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.541"
$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("E42").Value = "  -0.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.65"
$ws.Range("E43").Value = "  -3.55%  "

$ws.Range("E44").Value = "  +6.94%  "

$ws.Range("E45").Value = "  +1.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.37"
$ws.Range("E46").Value = "  +0.01%  "

$ws.Range("D47").Value = "1.731.91"
$ws.Range("E47").Value = "  +1.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.69"
$ws.Range("E48").Value = "  +2.08%  "

$ws.Range("E49").Value = "  +1.11%  "

$ws.Range("E50").Value = "  -1.45%  "

$ws.Range("E51").Value = "  -0.16%  "
